$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.878.87"
$ws.Range("E2").Value = "  -1.03%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.528.79"
$ws.Range("E3").Value = "  -0.89%  "

$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.21"
$ws.Range("E5").Value = "  -1.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.77"
$ws.Range("E6").Value = "  -0.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.527.02"
$ws.Range("E7").Value = "  -0.88%  "

$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("E10").Value = "  -1.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.79"
$ws.Range("E11").Value = "  -2.46%  "

$ws.Range("E12").Value = "  -2.37%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.127.75"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("E14").Value = "  -3.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.63"
$ws.Range("E15").Value = "  -4.64%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.526.21"
$ws.Range("E16").Value = "  -0.67%  "

$ws.Range("E17").Value = "  +1.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "65.848.60"
$ws.Range("E18").Value = "  -1.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.90"
$ws.Range("E19").Value = "  -5.58%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  -0.34%  "

$ws.Range("E21").Value = "  -4.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "414.00"
$ws.Range("E22").Value = "  -4.03%  "

$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.28"
$ws.Range("E24").Value = "  -2.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.670.15"

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("E27").Value = "  -2.83%  "

$ws.Range("E28").Value = "  -2.58%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.79"
$ws.Range("E29").Value = "  -2.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.90"
$ws.Range("E30").Value = "  -2.85%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  -0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.525.21"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.154"
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.33"
$ws.Range("E34").Value = "  -3.90%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.52"
$ws.Range("E36").Value = "  -4.04%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.28"
$ws.Range("E37").Value = "  -11.98%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "175.64"
$ws.Range("E38").Value = "  +0.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.26"
$ws.Range("E39").Value = "  -6.11%  "

$ws.Range("E40").Value = "  -8.29%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0819"
$ws.Range("E41").Value = "  -3.41%  "

$ws.Range("E42").Value = "  -2.50%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.860"
$ws.Range("E43").Value = "  -3.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.34"
$ws.Range("E44").Value = "  -1.67%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.78"
$ws.Range("E45").Value = "  -7.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.999"
$ws.Range("E46").Value = "  -0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("E47").Value = "  -4.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.06"
$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.86"
$ws.Range("E49").Value = "  -2.71%  "

$ws.Range("E50").Value = "  -7.38%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.00"
$ws.Range("E51").Value = "  -8.23%  "
